$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New localization rows appended at the bottom of the table (rows 119-122)
$data = @(
    @("lang_phone_number_collapse", "SDT", "Phone.No"),
    @("lang_phone_number_expand", "Số Điện Thoại", "Phone Number"),
    @("lang_role", "Vai Trò", "Role"),
    @("lang_guardians_list", "Danh Sách Giám Hộ", "Guardians List")
)

$startRow = 119
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Match the existing style used for "section header" key cells in column A
# (copy just the cell format from a row that already uses it, then paste
# formats only so the newly-written value/type is left untouched)
$ws.Range("A111").Copy()
$ws.Range("A122").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C122").Select()
